# Power BI Updates 24092023
#
# The source workbook is a plain data range (A1:F78) on a sheet named
# "district". The target edit:
#   1. Converts that range into a native Excel Table ("ListObject") named
#      "district", with autofilter + TableStyleMedium2 styling (xl/tables/table1.xml
#      + <tableParts> back-reference on the worksheet).
#   2. Renames the worksheet tab itself from "district" to "Sheet1".
#   3. Moves the active selection from D14 to D8.
#   4. Adjusts a handful of column widths (A, C, D, F).
#
# (Misc low-level version-stamp bytes in the diff -- fileVersion/rupBuild,
#  the x15ac:absPath scratch-dir string, xr:revisionPtr GUIDs, bookViews
#  xWindow/yWindow, and the x14ac:knownFonts / per-row x14ac:dyDescent
#  attributes -- are artifacts the real Excel.exe build stamps into the
#  OOXML on save and are not reachable through the Excel object model /
#  COM surface, so they're intentionally left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Turn A1:F78 into a real Table (ListObject), header row already present.
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:F78"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "district"

# 2) Rename the sheet tab.
$ws.Name = "Sheet1"

# 3) Column width tweaks (A/C/D/F). ColumnWidth is in "characters"; Excel's
#    column model quantizes the exported <col width> to an integer pixel
#    grid, so these character-width inputs are chosen to land on (or as
#    close as representable to) the target exported widths of 13.6640625,
#    12.5, 13.5 and 9.5 respectively.
$ws.Columns.Item(1).ColumnWidth = 12.75
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 8.666666666666666

# 4) Move the selection from D14 to D8.
$ws.Range("D8").Select()
